# Auto-generated script applying scheduled market-data refresh
# to the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# For each affected row: updates current market price columns
# (H/I/J/K/L), recomputed profit columns (M/N), clearing cells that
# no longer have a value and adding cells that newly do.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3000
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H74").Value = 5255
$ws.Range("I74").Value = 5466.923
$ws.Range("J74").Value = 2500
$ws.Range("K74").Value = 5466.923
$ws.Range("L74").Value = 2500
$ws.Range("M74").Value = -4530.923
$ws.Range("N74").Value = -4372
$ws.Range("H77").Value = 5255
$ws.Range("I77").Value = 5466.923
$ws.Range("J77").Value = 2500
$ws.Range("K77").Value = 27334.615
$ws.Range("L77").Value = 12500
$ws.Range("M77").Value = -22654.615
$ws.Range("N77").Value = -21860
$ws.Range("H86").Value = 4374.857
$ws.Range("I86").Value = 4675
$ws.Range("J86").Value = 3974.6667
$ws.Range("K86").Value = 4675
$ws.Range("L86").Value = 3974.6667
$ws.Range("M86").Value = -3552
$ws.Range("N86").Value = -6220.6667
$ws.Range("H89").Value = 4374.857
$ws.Range("I89").Value = 4675
$ws.Range("J89").Value = 3974.6667
$ws.Range("K89").Value = 23375
$ws.Range("L89").Value = 19873.3335
$ws.Range("M89").Value = -17759
$ws.Range("N89").Value = -31105.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 21390
$ws.Range("J44").Value = 21390
$ws.Range("L44").Value = 21390
$ws.Range("N44").Value = -22366
$ws.Range("H45").Value = 5463.3335
$ws.Range("I45").Value = 4990
$ws.Range("K45").Value = 4990
$ws.Range("M45").Value = -4613
$ws.Range("H61").Value = 5004723
$ws.Range("I61").Value = 4733.0586
$ws.Range("K61").Value = 4733.0586
$ws.Range("M61").Value = -4521.0586
$ws.Range("H74").Value = 1467845.5
$ws.Range("I74").Value = 1988621
$ws.Range("K74").Value = 1988621
$ws.Range("M74").Value = -1987747
$ws.Range("H77").Value = 1467845.5
$ws.Range("I77").Value = 1988621
$ws.Range("K77").Value = 9943105
$ws.Range("M77").Value = -9938737
$ws.Range("H132").Value = 1533986.8
$ws.Range("I132").Value = 1817350.9
$ws.Range("J132").Value = 3819.8
$ws.Range("K132").Value = 5452052.699999999
$ws.Range("L132").Value = 11459.4
$ws.Range("M132").Value = -5449522.699999999
$ws.Range("N132").Value = -16519.4
$ws.Range("H136").Value = 5004723
$ws.Range("I136").Value = 4733.0586
$ws.Range("K136").Value = 14199.1758
$ws.Range("M136").Value = -11649.1758

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 4002.5
$ws.Range("I10").Value = 4002.5
$ws.Range("K10").Value = 4002.5
$ws.Range("M10").Value = -3862.5
$ws.Range("H107").Value = 1987.5625
$ws.Range("J107").Value = 1989
$ws.Range("L107").Value = 1989
$ws.Range("N107").Value = -5829
$ws.Range("H134").Value = 2689472.2
$ws.Range("I134").Value = 1321.4263
$ws.Range("K134").Value = 3964.2789
$ws.Range("M134").Value = -1429.2789

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 68234.53
$ws.Range("I16").Value = 1424.4615
$ws.Range("K16").Value = 1424.4615
$ws.Range("M16").Value = -1137.4615
$ws.Range("H31").Value = 1391829.1
$ws.Range("I31").Value = 1504369.2
$ws.Range("J31").Value = 3833.3333
$ws.Range("K31").Value = 1504369.2
$ws.Range("L31").Value = 3833.3333
$ws.Range("M31").Value = -1504074.2
$ws.Range("N31").Value = -4423.3333
$ws.Range("H34").Value = 1391829.1
$ws.Range("I34").Value = 1504369.2
$ws.Range("J34").Value = 3833.3333
$ws.Range("K34").Value = 1504369.2
$ws.Range("L34").Value = 3833.3333
$ws.Range("M34").Value = -1504167.2
$ws.Range("N34").Value = -4237.3333
$ws.Range("H106").Value = 44411.5
$ws.Range("J106").Value = 44411.5
$ws.Range("L106").Value = 44411.5
$ws.Range("N106").Value = -46935.5
$ws.Range("H113").Value = 68234.53
$ws.Range("I113").Value = 1424.4615
$ws.Range("K113").Value = 1424.4615
$ws.Range("M113").Value = 745.5385000000001
$ws.Range("H132").Value = 2443.1428
$ws.Range("I132").Value = 1836.579
$ws.Range("J132").Value = 3723.6667
$ws.Range("K132").Value = 5509.737
$ws.Range("L132").Value = 11171.0001
$ws.Range("M132").Value = -2979.737
$ws.Range("N132").Value = -16231.0001
$ws.Range("H134").Value = 3666.7097
$ws.Range("I134").Value = 3421.7144
$ws.Range("J134").Value = 4181.2
$ws.Range("K134").Value = 10265.1432
$ws.Range("L134").Value = 12543.6
$ws.Range("M134").Value = -7730.143199999999
$ws.Range("N134").Value = -17613.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 1301.6666
$ws.Range("I60").Value = 577.5
$ws.Range("K60").Value = 1732.5
$ws.Range("M60").Value = -1481.5
$ws.Range("H64").Value = 2922.25
$ws.Range("I64").Value = 1566.3334
$ws.Range("K64").Value = 4699.0002
$ws.Range("M64").Value = -4429.0002
$ws.Range("H67").Value = 2922.25
$ws.Range("I67").Value = 1566.3334
$ws.Range("K67").Value = 4699.0002
$ws.Range("M67").Value = -3763.0002
$ws.Range("H92").Value = 844.8276
$ws.Range("J92").Value = 1857.1428
$ws.Range("L92").Value = 5571.428400000001
$ws.Range("N92").Value = -8067.428400000001
$ws.Range("H103").Value = 116
$ws.Range("I103").Value = 123.4
$ws.Range("J103").Value = 79
$ws.Range("K103").Value = 370.2
$ws.Range("L103").Value = 237
$ws.Range("M103").Value = 508.8
$ws.Range("N103").Value = -1995
$ws.Range("H107").Value = 554.625
$ws.Range("I107").Value = 304.08334
$ws.Range("J107").Value = 1306.25
$ws.Range("K107").Value = 912.2500200000001
$ws.Range("L107").Value = 3918.75
$ws.Range("M107").Value = 1007.74998
$ws.Range("N107").Value = -7758.75
$ws.Range("H137").Value = 3314.0908
$ws.Range("J137").Value = 4648.75
$ws.Range("L137").Value = 13946.25
$ws.Range("N137").Value = -24146.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 40173
$ws.Range("J34").Value = 40173
$ws.Range("L34").Value = 40173
$ws.Range("N34").Value = -40709
$ws.Range("H58").Value = 14966.667
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 14966.667
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 14966.667
$ws.Range("N58").Value = -15520.667
$ws.Range("M58").ClearContents()
$ws.Range("H76").Value = 40173
$ws.Range("J76").Value = 40173
$ws.Range("L76").Value = 40173
$ws.Range("N76").Value = -40803
$ws.Range("H79").Value = 40173
$ws.Range("J79").Value = 40173
$ws.Range("L79").Value = 40173
$ws.Range("N79").Value = -42357
$ws.Range("H107").Value = 690.8261
$ws.Range("I107").Value = 619
$ws.Range("K107").Value = 619
$ws.Range("M107").Value = 1301
$ws.Range("H132").Value = 5352.2104
$ws.Range("I132").Value = 5787.794
$ws.Range("K132").Value = 17363.382
$ws.Range("M132").Value = -14833.382

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3485.8125
$ws.Range("I7").Value = 2513.2222
$ws.Range("K7").Value = 2513.2222
$ws.Range("M7").Value = -2401.2222
$ws.Range("H20").Value = 27833.334
$ws.Range("I20").Value = 11750
$ws.Range("K20").Value = 11750
$ws.Range("M20").Value = -11524
$ws.Range("H40").Value = 2435
$ws.Range("I40").Value = 2438.4666
$ws.Range("J40").Value = 2429.8
$ws.Range("K40").Value = 2438.4666
$ws.Range("L40").Value = 2429.8
$ws.Range("M40").Value = -2302.4666
$ws.Range("N40").Value = -2701.8
$ws.Range("H43").Value = 313400
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H68").Value = 202057.8
$ws.Range("I68").Value = 251572.25
$ws.Range("J68").Value = 4000
$ws.Range("K68").Value = 251572.25
$ws.Range("L68").Value = 4000
$ws.Range("M68").Value = -250823.25
$ws.Range("N68").Value = -5498
$ws.Range("H71").Value = 202057.8
$ws.Range("I71").Value = 251572.25
$ws.Range("J71").Value = 4000
$ws.Range("K71").Value = 1257861.25
$ws.Range("L71").Value = 20000
$ws.Range("M71").Value = -1254117.25
$ws.Range("N71").Value = -27488
$ws.Range("H126").Value = 3485.8125
$ws.Range("I126").Value = 2513.2222
$ws.Range("K126").Value = 7539.6666
$ws.Range("M126").Value = -5069.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 136353.78
$ws.Range("I24").Value = 217503.6
$ws.Range("J24").Value = 34916.5
$ws.Range("K24").Value = 217503.6
$ws.Range("L24").Value = 34916.5
$ws.Range("M24").Value = -217273.6
$ws.Range("N24").Value = -35376.5
$ws.Range("H25").Value = 34166.668
$ws.Range("I25").Value = 30000
$ws.Range("K25").Value = 30000
$ws.Range("M25").Value = -29707
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H126").Value = 2715.9565
$ws.Range("I126").Value = 2191.2
$ws.Range("K126").Value = 6573.599999999999
$ws.Range("M126").Value = -4103.599999999999
$ws.Range("H132").Value = 3705091
$ws.Range("I132").Value = 4386898.5
$ws.Range("J132").Value = 3849.5715
$ws.Range("K132").Value = 13160695.5
$ws.Range("L132").Value = 11548.7145
$ws.Range("M132").Value = -13158165.5
$ws.Range("N132").Value = -16608.7145
$ws.Range("H136").Value = 4127700.5
$ws.Range("J136").Value = 33335408
$ws.Range("L136").Value = 100006224
$ws.Range("N136").Value = -100011324
